$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date from 2023-09-02 to 2023-09-03 for rows 2-7
$newDate = Get-Date -Year 2023 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
